$wb = $excel.ActiveWorkbook

# 1. Update Metadata timestamp
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 11:19 AM"

# 2. Update Industry Analysis sheet column F (1 Year) values
$wsIA = $wb.Worksheets.Item("Industry Analysis")
$wsIA.Range("F2").Value = 18.476
$wsIA.Range("F3").Value = -7.7404
$wsIA.Range("F4").Value = 30.7972
$wsIA.Range("F5").Value = -50.2266
$wsIA.Range("F6").Value = 61.9649
$wsIA.Range("F7").Value = -9.1713
$wsIA.Range("F8").Value = -3.556
$wsIA.Range("F9").Value = 38.3509
$wsIA.Range("F10").Value = -6.2497
$wsIA.Range("F11").Value = 52.6723
$wsIA.Range("F12").Value = -6.932
$wsIA.Range("F13").Value = 17.5662
$wsIA.Range("F14").Value = -35.5106
$wsIA.Range("F15").Value = 0.6286
$wsIA.Range("F16").Value = -3.1514
$wsIA.Range("F17").Value = -20.6354
$wsIA.Range("F18").Value = -0.0175
$wsIA.Range("F19").Value = -26.9255
$wsIA.Range("F20").Value = 44.703
$wsIA.Range("F21").Value = 10.0506
$wsIA.Range("F22").Value = 84.6016
$wsIA.Range("F23").Value = -54.4868
$wsIA.Range("F24").Value = -12.8122
$wsIA.Range("F25").Value = -9.182700000000001
$wsIA.Range("F26").Value = 5.9529
$wsIA.Range("F27").Value = -33.2998
$wsIA.Range("F28").Value = -20.4441
$wsIA.Range("F29").Value = -17.1514
$wsIA.Range("F30").Value = 24.527
$wsIA.Range("F31").Value = 57.6193
$wsIA.Range("F32").Value = -1.527
$wsIA.Range("F33").Value = -5.2378
$wsIA.Range("F34").Value = 27.4054
$wsIA.Range("F35").Value = 6.7961
$wsIA.Range("F36").Value = -5.6683
$wsIA.Range("F37").Value = 1.4178
$wsIA.Range("F38").Value = -22.4272
$wsIA.Range("F39").Value = 12.3741
$wsIA.Range("F40").Value = -5.138
$wsIA.Range("F41").Value = -0.1825
$wsIA.Range("F42").Value = 23.2483
$wsIA.Range("F43").Value = 14.456
$wsIA.Range("F44").Value = -11.1739
$wsIA.Range("F45").Value = 27.112
$wsIA.Range("F46").Value = -5.6252
$wsIA.Range("F47").Value = -36.5148
$wsIA.Range("F48").Value = -27.8397
$wsIA.Range("F49").Value = -25.4424
$wsIA.Range("F50").Value = -49.1173
$wsIA.Range("F51").Value = -51.065
$wsIA.Range("F52").Value = -35.4517
$wsIA.Range("F53").Value = -11.9879
$wsIA.Range("F54").Value = -3.0992
$wsIA.Range("F55").Value = -15.3441
$wsIA.Range("F56").Value = -25.937
$wsIA.Range("F57").Value = -29.1486
$wsIA.Range("F58").Value = -6.4093
$wsIA.Range("F59").Value = -23.3046
$wsIA.Range("F60").Value = -11.2657
$wsIA.Range("F61").Value = -9.777699999999999
$wsIA.Range("F62").Value = -16.0561
$wsIA.Range("F63").Value = -9.932499999999999
$wsIA.Range("F64").Value = 51.8767
$wsIA.Range("F65").Value = -43.5191
$wsIA.Range("F66").Value = 13.7315
$wsIA.Range("F67").Value = 12.6111
$wsIA.Range("F68").Value = 31.7532
$wsIA.Range("F69").Value = -19.9577
$wsIA.Range("F70").Value = -12.9642
$wsIA.Range("F71").Value = 13.2432
$wsIA.Range("F72").Value = 2.8232
$wsIA.Range("F73").Value = -9.179
$wsIA.Range("F74").Value = -14.2931
$wsIA.Range("F75").Value = 28.3699
$wsIA.Range("F76").Value = 45.5868

# 3. Update Stock List sheet: insert a new top row (CAPTRU-RE1) and drop the
#    last row (TRAVELFOOD), shifting every other row down by one.
$wsSL = $wb.Worksheets.Item("Stock List")
$iconVal = $wsSL.Range("A3").Value()
$wsSL.Rows.Item(2).Insert()
$wsSL.Rows.Item(77).Delete()
$wsSL.Range("A2:H2").ClearFormats()
$wsSL.Range("A2").Value = $iconVal
$wsSL.Range("B2").Value = "CAPTRU-RE1"
$wsSL.Range("C2").Value = "CAPTRU-RE1"
$wsSL.Range("D2").Value = 5.67
$wsSL.Range("E2").Value = -11.9565
$wsSL.Range("F2").Value = "N/A"
$wsSL.Range("G2").Value = "N/A"
$wsSL.Range("H2").Value = 0
